$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update handoff/handback datetimes for the 512d7bce... row (row 2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-09 03:55:08"
$wsZhCn.Range("G2").Value = "2016-01-09 03:55:51"

# "de-de" sheet: update handoff/handback datetimes for the 512d7bce... row (row 2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-09 03:55:17"
$wsDeDe.Range("G2").Value = "2016-01-09 03:56:07"
